$d = $word.ActiveDocument

# The document has an empty paragraph immediately after the paragraph
# that contains the "microservices_basics" GitHub hyperlink. That
# stray empty paragraph (and its paragraph mark) should be removed so
# the hyperlink paragraph is directly followed by the next
# ("Logging-service...") paragraph.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "microservices_basics") {
        if ($i + 1 -le $count) {
            $empty = $d.Paragraphs.Item($i + 1)
            if ($empty.Range.Text.Trim() -eq "") {
                $empty.Range.Delete()
            }
        }
        break
    }
}
